# Trade #19 closed at 2026-02-17 08:21:47 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.41   # Current Capital
$summary.Range("B4").Value = -0.59     # Total P&L $
$summary.Range("B5").Value = -0.62     # Total P&L %
$summary.Range("B6").Value = 19        # Total Trades
$summary.Range("B8").Value = 10        # Losing Trades
$summary.Range("B9").Value = 26.32     # Win Rate %

# ---- Strategy Status sheet (MarketMaking row) ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.41      # Capital
$status.Range("D4").Value = 19         # Trades
$status.Range("E4").Value = -0.59      # P&L $
$status.Range("F4").Value = -0.59      # P&L %
$status.Range("G4").Value = 26.32      # Win Rate %

# ---- Trade #19 row (row 20) closes out on both "All Trades" and "MarketMaking" sheets ----
$sheetsWithTrade = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetsWithTrade) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("G20").Value = 0.39            # Exit Price
    $ws.Range("H20").Value = "CLOSED"        # Status
    $ws.Range("I20").Value = -60.2041        # P&L %
    $ws.Range("J20").Value = -0.59           # P&L $
    $ws.Range("K20").Value = 99.41           # Capital After
    $ws.Range("P20").Value = "early_exit"    # Exit Reason
    $ws.Range("Q20").Value = 2.19            # Duration (min)
}
